$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")

# Move the "Outstanding" column (O) one column to the right (to P), leaving O
# empty, while keeping the existing column width/definitions ("cols") untouched.
$null = $wsRepay.Range("O1:O14").Copy()
$null = $wsRepay.Range("P1:P14").PasteSpecial(-4122)   # xlPasteFormats
$null = $wsRepay.Range("O1:O14").Copy()
$null = $wsRepay.Range("P1:P14").PasteSpecial(-4163)   # xlPasteValues
$wsRepay.Range("O1:O14").Clear()
$excel.CutCopyMode = 0

# Correct a couple of repayment schedule figures for installment #5 (row 7)
$wsRepay.Range("H7").Value = 65.2
$wsRepay.Range("K7").Value = 898.53
$wsRepay.Range("P7").Value = 898.53

# Make "Repayment Schedule" the active sheet/tab, with G9 selected
$null = $wsRepay.Activate()
$null = $wsRepay.Range("G9").Select()
